# Update the ClassName list on the active worksheet.
# The final column (A1:A25) is the header "ClassName" followed by the
# full, alphabetically-sorted, de-duplicated list of class names
# (replacing "7N-Monitor Mesa Pasture" with "Alakli Allotment" and
# inserting several new entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "ClassName",
    "Alakli Allotment",
    "Alkali Allotment",
    "Big Grayback Allotment",
    "Big Greyback Allotment",
    "Bigelow Lakes Pasture",
    "Butler Butte Allotment",
    "Crater Allotment",
    "Elliot Creek Allotment",
    "Fish Lake Allotment",
    "Hershberger Allotment",
    "High Cascade Ranger District",
    "High Cascades Ranger District",
    "Highcascades Ranger District",
    "Local",
    "Moist Meadow Pasture",
    "Rogue River National Forest",
    "Rogue River-Siskiyou National Forest",
    "Rouge River National Forest",
    "Siskiyou Mountains Ranger District",
    "Tiller Ranger District",
    "Umpqua National Forest",
    "Unknown Pasture",
    "Upper Big Applegate Allotment",
    "Woodruff Allotment"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
